$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41/42: dogwifhat / Maker swap places in the ranking ---
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'

# --- Price (column D) updates ---
$ws.Range("D2").Value = '66.051.46'
$ws.Range("D3").Value = '3.236.70'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.30'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.49'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D8").Value = '3.226.41'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.10'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.486'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.64'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000233'
$ws.Range("D15").Value = '3.757.09'
$ws.Range("D16").Value = '66.154.96'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '546.56'
$ws.Range("D18").Value = '3.244.29'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.09'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.48'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.742'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.82'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.43'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.07'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.32'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.97'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.66'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.75'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '568.47'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.65'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.34'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.15'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0858'
$ws.Range("D41").Value = '3.190.00'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.95'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.58'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.282'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.30'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.39'
$ws.Range("D48").Value = '0.0₃0556'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.63'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.112'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.20'

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("E3").Value = '  +7.10%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  +5.17%  '
$ws.Range("E6").Value = '  +8.98%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +7.04%  '
$ws.Range("E9").Value = '  +6.31%  '
$ws.Range("E10").Value = '  +12.43%  '
$ws.Range("E11").Value = '  +7.38%  '
$ws.Range("E12").Value = '  +6.08%  '
$ws.Range("E13").Value = '  +4.67%  '
$ws.Range("E14").Value = '  +7.03%  '
$ws.Range("E15").Value = '  +7.15%  '
$ws.Range("E16").Value = '  +3.17%  '
$ws.Range("E17").Value = '  +14.18%  '
$ws.Range("E18").Value = '  +7.24%  '
$ws.Range("E19").Value = '  +3.07%  '
$ws.Range("E20").Value = '  +7.25%  '
$ws.Range("E22").Value = '  +9.36%  '
$ws.Range("E23").Value = '  +11.37%  '
$ws.Range("E24").Value = '  +8.04%  '
$ws.Range("E25").Value = '  +4.18%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +20.43%  '
$ws.Range("E28").Value = '  +10.11%  '
$ws.Range("E29").Value = '  +7.55%  '
$ws.Range("E30").Value = '  +7.98%  '
$ws.Range("E31").Value = '  +6.71%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  +6.45%  '
$ws.Range("E34").Value = '  +10.14%  '
$ws.Range("E36").Value = '  +7.75%  '
$ws.Range("E37").Value = '  +5.22%  '
$ws.Range("E38").Value = '  +13.62%  '
$ws.Range("E39").Value = '  +8.80%  '
$ws.Range("E40").Value = '  +7.27%  '
$ws.Range("E41").Value = '  +11.25%  '
$ws.Range("E42").Value = '  +9.73%  '
$ws.Range("E43").Value = '  +4.53%  '
$ws.Range("E44").Value = '  +17.97%  '
$ws.Range("E45").Value = '  +12.11%  '
$ws.Range("E46").Value = '  +7.21%  '
$ws.Range("E48").Value = '  +6.53%  '
$ws.Range("E49").Value = '  +4.92%  '
$ws.Range("E50").Value = '  +4.54%  '
$ws.Range("E51").Value = '  +9.84%  '
